# The "MANGO" entry (row 5, A5) is being removed from the fruit list.
# Select the cell first (mirrors the interactive click-then-delete flow)
# and then clear its contents - this leaves row 5 empty while rows 6/7
# keep their original row numbers, and drops the now-unused "MANGO"
# shared string from the shared strings table on save.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("A5").Select()
$ws.Range("A5").ClearContents()
